$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the regression coefficient values (A Lag / C Lag rows) to reflect
# the new monthly-diff / 12-lag test-statistic results.
$ws.Range("B2").Value = "-0.372***"
$ws.Range("C2").Value = "0.01*"
$ws.Range("B3").Value = "-3.464***"
$ws.Range("C3").Value = "-0.808***"
